# For each year block (12 monthly rows: Jan..Dec), move the last three
# months (Oct, Nov, Dec) to the front of the block, pushing Jan..Sep down
# by three rows. This reproduces the row re-ordering seen in the diff for
# every year (2014, 2015, 2016, 2017), while leaving the header row (1)
# and all cell formatting/styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yearStarts = @(2, 14, 26, 38)   # first data row of each 12-month block
$scratchRow = 51                 # just past the last used row (49); cleared afterwards

foreach ($start in $yearStarts) {
    $end = $start + 11

    # Stage the whole 12-month block (columns A:I) in a scratch area so the
    # in-place rotation below doesn't clobber source data before it's read.
    $ws.Range("A" + $start + ":I" + $end).Copy()
    $ws.Range("A" + $scratchRow).PasteSpecial()

    # Oct/Nov/Dec (the last three staged rows) -> top of the block.
    $ws.Range("A" + ($scratchRow + 9) + ":I" + ($scratchRow + 11)).Copy()
    $ws.Range("A" + $start).PasteSpecial()

    # Jan..Sep (the first nine staged rows) -> rows 4..12 of the block.
    $ws.Range("A" + $scratchRow + ":I" + ($scratchRow + 8)).Copy()
    $ws.Range("A" + ($start + 3)).PasteSpecial()

    # Remove the scratch copy.
    $ws.Range("A" + $scratchRow + ":I" + ($scratchRow + 11)).Clear()
}
